{"js": "// Add a \"Meta description\" paragraph right after the title (H1), and move\n// the old trailing \"Play Book of Spells...\" / \"Discover the magic...\"\n// paragraphs into that new meta-description line, replacing the trailing\n// italic paragraph's text with the new image-generation prompt.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Find the first (title/H1) paragraph and insert the new meta\n//    description paragraph right after it.\nconst titlePara = paragraphs.items[0];\n\nconst metaText =\n  \"Meta description: Discover the magic of Book of Spells and play for free. \" +\n  \"Read our review on this 5-reel, 9-payline slot game, including the pros and cons and sound design.\";\n\n// Insert an empty paragraph first, give it the body (\"Normal\") style, then\n// fill in its text - this keeps the paragraph's formatting clean (matching\n// the rest of the document's body paragraphs) instead of inheriting the\n// Heading1 style of the title paragraph it was split off from.\nconst metaPara = titlePara.insertParagraph(\"\", Word.InsertLocation.after);\nmetaPara.style = \"Normal\";\nawait context.sync();\n\nmetaPara.insertText(metaText, Word.InsertLocation.start);\nawait context.sync();\n\n// Bold only the \"Meta description\" label (not the trailing colon/text).\nconst labelRange = metaPara.search(\"Meta description\", { matchCase: true });\nlabelRange.load(\"items\");\nawait context.sync();\nlabelRange.items[0].font.bold = true;\nawait context.sync();\n\n// 2) Drop the old duplicate bold title paragraph and replace the old\n//    italic summary paragraph's text with the new \"Prompt: ...\" text,\n//    keeping its existing (italic) formatting.\nconst refreshed = body.paragraphs;\nrefreshed.load(\"items/text\");\nawait context.sync();\n\nconst items = refreshed.items;\nconst lastIndex = items.length - 1;\nconst boldTitlePara = items[lastIndex - 1];\nconst italicSummaryPara = items[lastIndex];\n\nboldTitlePara.delete();\nawait context.sync();\n\nconst promptText =\n  'Prompt: Create a cartoon image featuring a happy Maya warrior with glasses to fit the theme of the game \"Book of Spells\". ' +\n  \"The Maya warrior should be standing in front of a cauldron with spell books and magic objects surrounding him. \" +\n  \"The background should be an enchanted forest with stars and sparkles. The image should be colorful and convey the theme of the game. \" +\n  \"Use bright colors for the Maya warrior's clothing and accessories, and make the cauldron and spell books stand out. \" +\n  \"The image should be eye-catching and visually striking to attract potential players to the game.\";\n\nitalicSummaryPara.insertText(promptText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add a \"Meta description\" paragraph right after the title (H1), and move\n# the old trailing \"Play Book of Spells...\" / \"Discover the magic...\"\n# paragraphs into that new meta-description line, replacing the trailing\n# italic paragraph's text with the new image-generation prompt.\n\n$d = $word.ActiveDocument\n\n# --- Part 1: insert the new \"Meta description\" paragraph after the title ---\n$titlePara = $d.Paragraphs(1)\n$insertPoint = $titlePara.Range\n$insertPoint.Collapse(0)          # wdCollapseEnd\n$insertPoint.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaPara.Range.Style = \"Normal\"\n$metaPara.Range.InsertBefore(\"Meta description: Discover the magic of Book of Spells and play for free. Read our review on this 5-reel, 9-payline slot game, including the pros and cons and sound design.\")\n\n# Bold only the \"Meta description\" label (not the trailing colon/text).\n$metaPara = $d.Paragraphs(2)\n$labelRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + 16)\n$labelRange.Bold = 1\n\n# --- Part 2: drop the old duplicate bold title paragraph at the end, and ---\n# --- replace the old italic summary paragraph's text with the new prompt ---\n$count = $d.Paragraphs.Count\n$oldTitlePara = $d.Paragraphs($count - 1)\n$oldTitlePara.Range.Delete()\n\n$count = $d.Paragraphs.Count\n$summaryPara = $d.Paragraphs($count)\n$summaryRange = $d.Range($summaryPara.Range.Start, $summaryPara.Range.End)\n$promptText = \"Prompt: Create a cartoon image featuring a happy Maya warrior with glasses to fit the theme of the game \"\"Book of Spells\"\". The Maya warrior should be standing in front of a cauldron with spell books and magic objects surrounding him. The background should be an enchanted forest with stars and sparkles. The image should be colorful and convey the theme of the game. Use bright colors for the Maya warrior's clothing and accessories, and make the cauldron and spell books stand out. The image should be eye-catching and visually striking to attract potential players to the game.\"\n$summaryRange.Text = $promptText\n"}
